$d = $word.ActiveDocument

# --- Change 1: insert a tab (as its own run) before the title text ---
$titlePara = $d.Paragraphs(1)
$titleStart = $titlePara.Range.Start
$insertPoint = $d.Range($titleStart, $titleStart)
$tabXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial Black" w:hAnsi="Arial Black"/></w:rPr><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($tabXml)
# InsertXML split the title paragraph in two (tab-only paragraph + original title
# paragraph); merge them back into a single paragraph by deleting the paragraph
# mark that now separates them.
$tabPara = $d.Paragraphs(1).Range
$mark = $d.Range($tabPara.End - 1, $tabPara.End)
$mark.Delete()

# --- Change 2: append new paragraphs at the end of the document ---
# First create a genuine new trailing empty paragraph after the current last
# paragraph (this keeps the existing last paragraph's identity/position intact).
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# Insert the new content into that fresh trailing paragraph.
$newLast = $d.Paragraphs($d.Paragraphs.Count)
$tailRange = $newLast.Range
$tailXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Ayer día 11/11/2025 descargamos una copia de Git</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> en la que sale una carpeta clonada</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> y la añadimos al visual donde los cambios que realicemos se puedan subir al Git a través de visual sin la necesidad de meterse en Git.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Esto yo lo añadí a través de un código dicho por el profesor en la propia terminal del Visual</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Vi también todo el rendimiento de mi página con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>LightHdus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> en la que por ejemplo pasar las imágenes de JPG a Web subiría el rendimiento.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">También se le añade in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>lazy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> a las imágenes para que me carguen mas.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRange.InsertXML($tailXml)

# InsertXML leaves one spurious empty paragraph at the very end (mirroring the
# tab-insertion behaviour above); merge it away the same way.
$secondLast = $d.Paragraphs($d.Paragraphs.Count - 1)
$secondLastRange = $secondLast.Range
$mark2 = $d.Range($secondLastRange.End - 1, $secondLastRange.End)
$mark2.Delete()
